$wb = $excel.ActiveWorkbook

# ---- Sheet ALC: 62 cell update(s) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 83338670  # H29: 100004810 -> 83338670
$ws.Cells.Item(40, 8).Value = 3688.6  # H40: 3876.7778 -> 3688.6
$ws.Cells.Item(40, 9).Value = 2497.75  # I40: 2998.5 -> 2497.75
$ws.Cells.Item(40, 10).Value = 4482.5  # J40: 4127.7144 -> 4482.5
$ws.Cells.Item(40, 11).Value = 2497.75  # K40: 2998.5 -> 2497.75
$ws.Cells.Item(40, 12).Value = 4482.5  # L40: 4127.7144 -> 4482.5
$ws.Cells.Item(40, 13).Value = -2322.75  # M40: -2823.5 -> -2322.75
$ws.Cells.Item(40, 14).Value = -4832.5  # N40: -4477.7144 -> -4832.5
$ws.Cells.Item(55, 8).Value = 290.4  # H55: 300.5 -> 290.4
$ws.Cells.Item(55, 9).Value = 287.5  # I55: 300 -> 287.5
$ws.Cells.Item(55, 11).Value = 287.5  # K55: 300 -> 287.5
$ws.Cells.Item(55, 13).Value = -73.5  # M55: -86 -> -73.5
$ws.Cells.Item(58, 8).Value = 20837970  # H58: 19235384 -> 20837970
$ws.Cells.Item(58, 9).Value = 22727786  # I58: 22727742 -> 22727786
$ws.Cells.Item(58, 10).Value = 49999  # J58: 27408 -> 49999
$ws.Cells.Item(58, 11).Value = 68183358  # K58: 68183226 -> 68183358
$ws.Cells.Item(58, 12).Value = 149997  # L58: 82224 -> 149997
$ws.Cells.Item(58, 13).Value = -68183208  # M58: -68183076 -> -68183208
$ws.Cells.Item(58, 14).Value = -150297  # N58: -82524 -> -150297
$ws.Cells.Item(74, 8).Value = 20238.834  # H74: 23486.6 -> 20238.834
$ws.Cells.Item(74, 9).Value = 20238.834  # I74: 23486.6 -> 20238.834
$ws.Cells.Item(74, 11).Value = 20238.834  # K74: 23486.6 -> 20238.834
$ws.Cells.Item(74, 13).Value = -19302.834  # M74: -22550.6 -> -19302.834
$ws.Cells.Item(77, 8).Value = 20238.834  # H77: 23486.6 -> 20238.834
$ws.Cells.Item(77, 9).Value = 20238.834  # I77: 23486.6 -> 20238.834
$ws.Cells.Item(77, 11).Value = 101194.17  # K77: 117433 -> 101194.17
$ws.Cells.Item(77, 13).Value = -96514.17  # M77: -112753 -> -96514.17
$ws.Cells.Item(86, 8).Value = 3087  # H86: 3079.6667 -> 3087
$ws.Cells.Item(86, 9).Value = 3020  # I86: 3036.6667 -> 3020
$ws.Cells.Item(86, 10).Value = 3167.4  # J86: 3122.6667 -> 3167.4
$ws.Cells.Item(86, 11).Value = 3020  # K86: 3036.6667 -> 3020
$ws.Cells.Item(86, 12).Value = 3167.4  # L86: 3122.6667 -> 3167.4
$ws.Cells.Item(86, 13).Value = -1897  # M86: -1913.6667 -> -1897
$ws.Cells.Item(86, 14).Value = -5413.4  # N86: -5368.6667 -> -5413.4
$ws.Cells.Item(88, 8).Value = 3280.318  # H88: 3289.8635 -> 3280.318
$ws.Cells.Item(88, 9).Value = 666.3333  # I88: 699.6667 -> 666.3333
$ws.Cells.Item(88, 10).Value = 3693.0527  # J88: 3698.842 -> 3693.0527
$ws.Cells.Item(88, 11).Value = 666.3333  # K88: 699.6667 -> 666.3333
$ws.Cells.Item(88, 12).Value = 3693.0527  # L88: 3698.842 -> 3693.0527
$ws.Cells.Item(88, 13).Value = -260.3333  # M88: -293.6667 -> -260.3333
$ws.Cells.Item(88, 14).Value = -4505.0527  # N88: -4510.842000000001 -> -4505.0527
$ws.Cells.Item(89, 8).Value = 3087  # H89: 3079.6667 -> 3087
$ws.Cells.Item(89, 9).Value = 3020  # I89: 3036.6667 -> 3020
$ws.Cells.Item(89, 10).Value = 3167.4  # J89: 3122.6667 -> 3167.4
$ws.Cells.Item(89, 11).Value = 15100  # K89: 15183.3335 -> 15100
$ws.Cells.Item(89, 12).Value = 15837  # L89: 15613.3335 -> 15837
$ws.Cells.Item(89, 13).Value = -9484  # M89: -9567.333500000001 -> -9484
$ws.Cells.Item(89, 14).Value = -27069  # N89: -26845.3335 -> -27069
$ws.Cells.Item(91, 8).Value = 3280.318  # H91: 3289.8635 -> 3280.318
$ws.Cells.Item(91, 9).Value = 666.3333  # I91: 699.6667 -> 666.3333
$ws.Cells.Item(91, 10).Value = 3693.0527  # J91: 3698.842 -> 3693.0527
$ws.Cells.Item(91, 11).Value = 666.3333  # K91: 699.6667 -> 666.3333
$ws.Cells.Item(91, 12).Value = 3693.0527  # L91: 3698.842 -> 3693.0527
$ws.Cells.Item(91, 13).Value = 737.6667  # M91: 704.3333 -> 737.6667
$ws.Cells.Item(91, 14).Value = -6501.0527  # N91: -6506.842000000001 -> -6501.0527
$ws.Cells.Item(141, 8).Value = 9223.666999999999  # H141: 5451.3335 -> 9223.666999999999
$ws.Cells.Item(141, 9).Value = 6118.3335  # I141: 5927 -> 6118.3335
$ws.Cells.Item(141, 10).Value = 10000  # J141: 4500 -> 10000
$ws.Cells.Item(141, 11).Value = 18355.0005  # K141: 17781 -> 18355.0005
$ws.Cells.Item(141, 12).Value = 30000  # L141: 13500 -> 30000
$ws.Cells.Item(141, 13).Value = -13175.0005  # M141: -12601 -> -13175.0005
$ws.Cells.Item(141, 14).Value = -40360  # N141: -23860 -> -40360

# ---- Sheet ARM: 16 cell update(s) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2243.8572  # H2: 2362.6155 -> 2243.8572
$ws.Cells.Item(2, 9).Value = 1273.5714  # I2: 1369.1666 -> 1273.5714
$ws.Cells.Item(2, 11).Value = 1273.5714  # K2: 1369.1666 -> 1273.5714
$ws.Cells.Item(2, 13).Value = -1160.5714  # M2: -1256.1666 -> -1160.5714
$ws.Cells.Item(36, 8).Value = 4466.8335  # H36: 4721.1665 -> 4466.8335
$ws.Cells.Item(36, 9).Value = 5450.25  # I36: 5165.4 -> 5450.25
$ws.Cells.Item(36, 11).Value = 5450.25  # K36: 5165.4 -> 5450.25
$ws.Cells.Item(36, 13).Value = -5104.25  # M36: -4819.4 -> -5104.25
$ws.Cells.Item(116, 8).Value = 2243.8572  # H116: 2362.6155 -> 2243.8572
$ws.Cells.Item(116, 9).Value = 1273.5714  # I116: 1369.1666 -> 1273.5714
$ws.Cells.Item(116, 11).Value = 1273.5714  # K116: 1369.1666 -> 1273.5714
$ws.Cells.Item(116, 13).Value = 1020.4286  # M116: 924.8334 -> 1020.4286
$ws.Cells.Item(132, 8).Value = 2080.6  # H132: 1866.7097 -> 2080.6
$ws.Cells.Item(132, 9).Value = 2080.6  # I132: 1866.7097 -> 2080.6
$ws.Cells.Item(132, 11).Value = 6241.799999999999  # K132: 5600.1291 -> 6241.799999999999
$ws.Cells.Item(132, 13).Value = -3711.799999999999  # M132: -3070.1291 -> -3711.799999999999

# ---- Sheet BSM: 4 cell update(s) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2243.8572  # H3: 2362.6155 -> 2243.8572
$ws.Cells.Item(3, 9).Value = 1273.5714  # I3: 1369.1666 -> 1273.5714
$ws.Cells.Item(3, 11).Value = 1273.5714  # K3: 1369.1666 -> 1273.5714
$ws.Cells.Item(3, 13).Value = -1159.5714  # M3: -1255.1666 -> -1159.5714

# ---- Sheet CRP: 26 cell update(s) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 2323.5  # H8: 2361 -> 2323.5
$ws.Cells.Item(8, 10).Value = 2323.5  # J8: 2361 -> 2323.5
$ws.Cells.Item(8, 12).Value = 2323.5  # L8: 2361 -> 2323.5
$ws.Cells.Item(8, 14).Value = -2603.5  # N8: -2641 -> -2603.5
$ws.Cells.Item(99, 8).Value = 2347.4736  # H99: 2394.5557 -> 2347.4736
$ws.Cells.Item(99, 9).Value = 2225.25  # I99: 2273.6 -> 2225.25
$ws.Cells.Item(99, 11).Value = 2225.25  # K99: 2273.6 -> 2225.25
$ws.Cells.Item(99, 13).Value = -727.25  # M99: -775.5999999999999 -> -727.25
$ws.Cells.Item(122, 8).Value = 3210.2542  # H122: 3291.5615 -> 3210.2542
$ws.Cells.Item(122, 9).Value = 821.875  # I122: 834.5625 -> 821.875
$ws.Cells.Item(122, 10).Value = 4098.9536  # J122: 4250.39 -> 4098.9536
$ws.Cells.Item(122, 11).Value = 2465.625  # K122: 2503.6875 -> 2465.625
$ws.Cells.Item(122, 12).Value = 12296.8608  # L122: 12751.17 -> 12296.8608
$ws.Cells.Item(122, 13).Value = -15.625  # M122: -53.6875 -> -15.625
$ws.Cells.Item(122, 14).Value = -17196.8608  # N122: -17651.17 -> -17196.8608
$ws.Cells.Item(126, 8).Value = 2347.4736  # H126: 2394.5557 -> 2347.4736
$ws.Cells.Item(126, 9).Value = 2225.25  # I126: 2273.6 -> 2225.25
$ws.Cells.Item(126, 11).Value = 6675.75  # K126: 6820.799999999999 -> 6675.75
$ws.Cells.Item(126, 13).Value = -4205.75  # M126: -4350.799999999999 -> -4205.75
$ws.Cells.Item(132, 8).Value = 2235.1924  # H132: 2126.5518 -> 2235.1924
$ws.Cells.Item(132, 9).Value = 2244.6  # I132: 2126.5518 -> 2244.6
$ws.Cells.Item(132, 10).Value = 2000  # J132: 0 -> 2000
$ws.Cells.Item(132, 11).Value = 6733.799999999999  # K132: 6379.655400000001 -> 6733.799999999999
$ws.Cells.Item(132, 12).Value = 6000  # L132: 0 -> 6000
$ws.Cells.Item(132, 13).Value = -4203.799999999999  # M132: -3849.655400000001 -> -4203.799999999999
$ws.Cells.Item(132, 14).Value = -11060  # N132: None -> -11060

# ---- Sheet CUL: 46 cell update(s) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 72.30768999999999  # H2: 67.21429000000001 -> 72.30768999999999
$ws.Cells.Item(2, 9).Value = 153  # I2: 115 -> 153
$ws.Cells.Item(2, 11).Value = 918  # K2: 690 -> 918
$ws.Cells.Item(2, 13).Value = -805  # M2: -577 -> -805
$ws.Cells.Item(3, 8).Value = 125003750  # H3: 71431420 -> 125003750
$ws.Cells.Item(3, 9).Value = 166669810  # I3: 100002280 -> 166669810
$ws.Cells.Item(3, 10).Value = 5555  # J3: 4277.5 -> 5555
$ws.Cells.Item(3, 11).Value = 500009430  # K3: 300006840 -> 500009430
$ws.Cells.Item(3, 12).Value = 16665  # L3: 12832.5 -> 16665
$ws.Cells.Item(3, 13).Value = -500009318  # M3: -300006728 -> -500009318
$ws.Cells.Item(3, 14).Value = -16889  # N3: -13056.5 -> -16889
$ws.Cells.Item(38, 8).Value = 669.1177  # H38: 634.2222 -> 669.1177
$ws.Cells.Item(38, 10).Value = 1019.7273  # J38: 938.1667 -> 1019.7273
$ws.Cells.Item(38, 12).Value = 3059.1819  # L38: 2814.5001 -> 3059.1819
$ws.Cells.Item(38, 14).Value = -3753.1819  # N38: -3508.5001 -> -3753.1819
$ws.Cells.Item(68, 8).Value = 1500  # H68: 1001 -> 1500
$ws.Cells.Item(68, 9).Value = 0  # I68: 502 -> 0
$ws.Cells.Item(68, 11).Value = 0  # K68: 1506 -> 0
$ws.Cells.Item(68, 13).ClearContents()  # M68 was -695
$ws.Cells.Item(71, 8).Value = 1500  # H71: 1001 -> 1500
$ws.Cells.Item(71, 9).Value = 0  # I71: 502 -> 0
$ws.Cells.Item(71, 11).Value = 0  # K71: 4518 -> 0
$ws.Cells.Item(71, 13).ClearContents()  # M71 was -462
$ws.Cells.Item(102, 8).Value = 15405.2  # H102: 9380.666999999999 -> 15405.2
$ws.Cells.Item(102, 10).Value = 18775  # J102: 10312.5 -> 18775
$ws.Cells.Item(102, 12).Value = 56325  # L102: 30937.5 -> 56325
$ws.Cells.Item(102, 14).Value = -61193  # N102: -35805.5 -> -61193
$ws.Cells.Item(105, 8).Value = 38958.5  # H105: 48315.332 -> 38958.5
$ws.Cells.Item(105, 10).Value = 38958.5  # J105: 48315.332 -> 38958.5
$ws.Cells.Item(105, 12).Value = 116875.5  # L105: 144945.996 -> 116875.5
$ws.Cells.Item(105, 14).Value = -122117.5  # N105: -150187.996 -> -122117.5
$ws.Cells.Item(112, 8).Value = 7500.8  # H112: 8849.75 -> 7500.8
$ws.Cells.Item(112, 10).Value = 8751  # J112: 10966.333 -> 8751
$ws.Cells.Item(112, 12).Value = 26253  # L112: 32898.999 -> 26253
$ws.Cells.Item(112, 14).Value = -28469  # N112: -35114.999 -> -28469
$ws.Cells.Item(113, 8).Value = 1188.2222  # H113: 1228.7059 -> 1188.2222
$ws.Cells.Item(113, 10).Value = 1833.1666  # J113: 2099.8 -> 1833.1666
$ws.Cells.Item(113, 12).Value = 5499.4998  # L113: 6299.400000000001 -> 5499.4998
$ws.Cells.Item(113, 14).Value = -9839.4998  # N113: -10639.4 -> -9839.4998
$ws.Cells.Item(132, 8).Value = 1457.4286  # H132: 2007.7307 -> 1457.4286
$ws.Cells.Item(132, 9).Value = 1457.4286  # I132: 1453.6 -> 1457.4286
$ws.Cells.Item(132, 10).Value = 0  # J132: 2763.3635 -> 0
$ws.Cells.Item(132, 11).Value = 13116.8574  # K132: 13082.4 -> 13116.8574
$ws.Cells.Item(132, 12).Value = 0  # L132: 24870.2715 -> 0
$ws.Cells.Item(132, 13).Value = -10586.8574  # M132: -10552.4 -> -10586.8574
$ws.Cells.Item(132, 14).ClearContents()  # N132 was -29930.2715

# ---- Sheet GSM: 26 cell update(s) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6075.7646  # H70: 6401.375 -> 6075.7646
$ws.Cells.Item(70, 9).Value = 6021.222  # I70: 6347.222 -> 6021.222
$ws.Cells.Item(70, 10).Value = 6137.125  # J70: 6471 -> 6137.125
$ws.Cells.Item(70, 11).Value = 6021.222  # K70: 6347.222 -> 6021.222
$ws.Cells.Item(70, 12).Value = 6137.125  # L70: 6471 -> 6137.125
$ws.Cells.Item(70, 13).Value = -5751.222  # M70: -6077.222 -> -5751.222
$ws.Cells.Item(70, 14).Value = -6677.125  # N70: -7011 -> -6677.125
$ws.Cells.Item(73, 8).Value = 6075.7646  # H73: 6401.375 -> 6075.7646
$ws.Cells.Item(73, 9).Value = 6021.222  # I73: 6347.222 -> 6021.222
$ws.Cells.Item(73, 10).Value = 6137.125  # J73: 6471 -> 6137.125
$ws.Cells.Item(73, 11).Value = 6021.222  # K73: 6347.222 -> 6021.222
$ws.Cells.Item(73, 12).Value = 6137.125  # L73: 6471 -> 6137.125
$ws.Cells.Item(73, 13).Value = -5085.222  # M73: -5411.222 -> -5085.222
$ws.Cells.Item(73, 14).Value = -8009.125  # N73: -8343 -> -8009.125
$ws.Cells.Item(97, 8).Value = 760.2  # H97: 782 -> 760.2
$ws.Cells.Item(97, 9).Value = 649.875  # I97: 677.4 -> 649.875
$ws.Cells.Item(97, 11).Value = 649.875  # K97: 677.4 -> 649.875
$ws.Cells.Item(97, 13).Value = -153.875  # M97: -181.4 -> -153.875
$ws.Cells.Item(126, 8).Value = 2436.5  # H126: 2570.625 -> 2436.5
$ws.Cells.Item(126, 9).Value = 2373.889  # I126: 2509.2856 -> 2373.889
$ws.Cells.Item(126, 11).Value = 7121.667  # K126: 7527.8568 -> 7121.667
$ws.Cells.Item(126, 13).Value = -4651.667  # M126: -5057.8568 -> -4651.667
$ws.Cells.Item(132, 8).Value = 3037  # H132: 3125.75 -> 3037
$ws.Cells.Item(132, 9).Value = 2838.75  # I132: 2933.682 -> 2838.75
$ws.Cells.Item(132, 11).Value = 8516.25  # K132: 8801.045999999998 -> 8516.25
$ws.Cells.Item(132, 13).Value = -5986.25  # M132: -6271.045999999998 -> -5986.25

# ---- Sheet LTW: 27 cell update(s) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2078  # H16: 1888.2 -> 2078
$ws.Cells.Item(16, 9).Value = 1203.5  # I16: 1057.2858 -> 1203.5
$ws.Cells.Item(16, 11).Value = 1203.5  # K16: 1057.2858 -> 1203.5
$ws.Cells.Item(16, 13).Value = -1033.5  # M16: -887.2858000000001 -> -1033.5
$ws.Cells.Item(22, 8).Value = 3422.2778  # H22: 3558.8823 -> 3422.2778
$ws.Cells.Item(22, 10).Value = 3279.3333  # J22: 3477.4546 -> 3279.3333
$ws.Cells.Item(22, 12).Value = 3279.3333  # L22: 3477.4546 -> 3279.3333
$ws.Cells.Item(22, 14).Value = -3869.3333  # N22: -4067.4546 -> -3869.3333
$ws.Cells.Item(27, 8).Value = 3422.2778  # H27: 3558.8823 -> 3422.2778
$ws.Cells.Item(27, 10).Value = 3279.3333  # J27: 3477.4546 -> 3279.3333
$ws.Cells.Item(27, 12).Value = 3279.3333  # L27: 3477.4546 -> 3279.3333
$ws.Cells.Item(27, 14).Value = -3493.3333  # N27: -3691.4546 -> -3493.3333
$ws.Cells.Item(32, 8).Value = 6044.8  # H32: 8603.429 -> 6044.8
$ws.Cells.Item(32, 9).Value = 6044.8  # I32: 8603.429 -> 6044.8
$ws.Cells.Item(32, 11).Value = 6044.8  # K32: 8603.429 -> 6044.8
$ws.Cells.Item(32, 13).Value = -5727.8  # M32: -8286.429 -> -5727.8
$ws.Cells.Item(43, 8).Value = 23161.514  # H43: 22799.375 -> 23161.514
$ws.Cells.Item(43, 9).Value = 22237.191  # I43: 21825.217 -> 22237.191
$ws.Cells.Item(43, 10).Value = 24374.688  # J43: 24117.354 -> 24374.688
$ws.Cells.Item(43, 11).Value = 22237.191  # K43: 21825.217 -> 22237.191
$ws.Cells.Item(43, 12).Value = 24374.688  # L43: 24117.354 -> 24374.688
$ws.Cells.Item(43, 13).Value = -22044.191  # M43: -21632.217 -> -22044.191
$ws.Cells.Item(43, 14).Value = -24760.688  # N43: -24503.354 -> -24760.688
$ws.Cells.Item(100, 8).Value = 8333.223  # H100: 8875 -> 8333.223
$ws.Cells.Item(100, 9).Value = 5666.3335  # I100: 6500 -> 5666.3335
$ws.Cells.Item(100, 11).Value = 5666.3335  # K100: 6500 -> 5666.3335
$ws.Cells.Item(100, 13).Value = -5125.3335  # M100: -5959 -> -5125.3335

# ---- Sheet WVR: 12 cell update(s) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1258.5652  # H113: 1262.6522 -> 1258.5652
$ws.Cells.Item(113, 9).Value = 836.2727  # I113: 844.8182 -> 836.2727
$ws.Cells.Item(113, 11).Value = 2508.8181  # K113: 2534.4546 -> 2508.8181
$ws.Cells.Item(113, 13).Value = -338.8181  # M113: -364.4546 -> -338.8181
$ws.Cells.Item(122, 8).Value = 4370.7144  # H122: 4385.143 -> 4370.7144
$ws.Cells.Item(122, 9).Value = 1556.0714  # I122: 1577.7142 -> 1556.0714
$ws.Cells.Item(122, 11).Value = 4668.2142  # K122: 4733.142599999999 -> 4668.2142
$ws.Cells.Item(122, 13).Value = -2218.2142  # M122: -2283.142599999999 -> -2218.2142
$ws.Cells.Item(132, 8).Value = 2256.111  # H132: 2288.375 -> 2256.111
$ws.Cells.Item(132, 9).Value = 2049.3333  # I132: 2059.6 -> 2049.3333
$ws.Cells.Item(132, 11).Value = 6147.999899999999  # K132: 6178.799999999999 -> 6147.999899999999
$ws.Cells.Item(132, 13).Value = -3617.999899999999  # M132: -3648.799999999999 -> -3617.999899999999
